# Fruta / hortaliza, semanal
# Insert a new weekly record at row 91 (shifting existing rows 91-182 down to 92-183)
# on the "Hortaliza, Terminal Hortofrutícola Agro Chillán - Repollo" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 91, pushing all data below it down by one row.
$ws.Rows.Item(91).Insert()

# Populate the newly inserted row 91 with the new weekly price record.
$ws.Cells.Item(91, 1).Value = 7
$ws.Cells.Item(91, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(91, 3).Value = "Ñuble"
$ws.Cells.Item(91, 4).Value = 44601
$ws.Cells.Item(91, 5).Value = 16
$ws.Cells.Item(91, 6).Value = 100112006
$ws.Cells.Item(91, 7).Value = "Repollo"
$ws.Cells.Item(91, 8).Value = "Crespo record"
$ws.Cells.Item(91, 9).Value = "Primera"
$ws.Cells.Item(91, 10).Value = 200
$ws.Cells.Item(91, 11).Value = 750
$ws.Cells.Item(91, 12).Value = 800
$ws.Cells.Item(91, 13).Value = 775
$ws.Cells.Item(91, 14).Value = "$/unidad"
$ws.Cells.Item(91, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(91, 16).Value = 775
$ws.Cells.Item(91, 17).Value = 1
$ws.Cells.Item(91, 18).Value = "Hortaliza"
